# Add a new "ODI Batting Extra" worksheet at the end of the workbook and
# populate it with the MATCH_CODE / BATTING_POSITION / NUM_4 / NUM_6 /
# PERCENT_RUNS_OF_TOTAL / MAN_OF_MATCH header row plus one data row.

$wb = $excel.ActiveWorkbook

# Grab a header cell from an existing sheet so the new header row can pick
# up the same bold / centered / bordered style already used by the other
# sheets (style index 1) instead of minting a near-duplicate style.
$styleSource = $wb.Worksheets.Item("ODI Bowling").Range("A1")

# Insert the new sheet after the current last sheet so it lands at the end.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "ODI Batting Extra"

# Header row.
$ws.Range("A1").Value = "MATCH_CODE"
$ws.Range("B1").Value = "BATTING_POSITION"
$ws.Range("C1").Value = "NUM_4"
$ws.Range("D1").Value = "NUM_6"
$ws.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$ws.Range("F1").Value = "MAN_OF_MATCH"

# Copy formatting only (xlPasteFormats = -4122) from the existing header
# cell onto the new header row.
$styleSource.Copy()
$ws.Range("A1:F1").PasteSpecial(-4122)

# Data row. MATCH_CODE is stored as text everywhere else in this workbook
# (Player Info!ID, ODI Batting!MATCH_CODE, ODI Bowling!MATCH_CODE are all
# text), so force the purely-numeric "4656" to stay text instead of being
# auto-converted to a number — a leading apostrophe is the normal Excel way
# to do that.
$ws.Range("A2").Value = "'4656"
$ws.Range("F2").Value = "NO"
